# Compcopy and order management record saving
#
# The template's sheet data used to start on row 2 (row 1 was a blank
# spacer row above the header), covering A2:Z6. This edit removes that
# leading blank row so the header + 4 data rows now occupy A1:Z5.
#
# Deleting the entire row 1 shifts every row's content up by one
# (values/styles are preserved as-is; only row/cell references move).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the blank first row - everything below shifts up one row.
$ws.Rows("1:1").Delete()

# The AutoFilter range doesn't auto-adjust from the row delete above,
# so reapply it over the new header/data extent (A1:Z5).
$ws.AutoFilterMode = $false
$ws.Range("A1:Z5").AutoFilter()

# Likewise, update the workbook-level hidden _FilterDatabase defined
# name (tied to the AutoFilter) so it points at the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$Z`$5"
    }
}

# Leave the cursor where it ended up after the edits.
$ws.Range("H15").Select()
